$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.242.35"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "2.284.08"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'323.89"
$ws.Range("E5").Value = "  +2.06%  "
$ws.Range("D6").Value = "'102.92"
$ws.Range("E6").Value = "  -1.90%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("D10").Value = "'39.82"
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("D11").Value = "'0.0907"
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").Value = "'8.33"
$ws.Range("E12").Value = "  -1.23%  "
$ws.Range("D13").Value = "'0.106"
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("D14").Value = "'0.971"
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("D15").Value = "'15.12"
$ws.Range("E15").Value = "  -2.20%  "
$ws.Range("D16").Value = "2.630.94"
$ws.Range("D17").Value = "2.287.01"
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("D18").Value = "42.228.54"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("D19").Value = "'7.36"
$ws.Range("E19").Value = "  -5.47%  "
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("D21").Value = "'13.03"
$ws.Range("E21").Value = "  +30.34%  "
$ws.Range("D22").Value = "'3.63"
$ws.Range("E22").Value = "  +1.93%  "
$ws.Range("D23").Value = "'73.06"
$ws.Range("E23").Value = "  -0.87%  "
$ws.Range("D24").Value = "'268.23"
$ws.Range("E24").Value = "  -6.28%  "
$ws.Range("E25").Value = "  -3.10%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("D28").Value = "'2.29"
$ws.Range("E28").Value = "  +2.86%  "
$ws.Range("D29").Value = "'22.45"
$ws.Range("E29").Value = "  -4.17%  "
$ws.Range("D30").Value = "'37.81"
$ws.Range("E30").Value = "  +6.49%  "
$ws.Range("D31").Value = "'164.01"
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("D32").Value = "'6.14"
$ws.Range("E32").Value = "  +3.61%  "
$ws.Range("D33").Value = "'0.0875"
$ws.Range("E33").Value = "  -0.97%  "
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("D35").Value = "'0.114"
$ws.Range("E35").Value = "  -2.23%  "
$ws.Range("E36").Value = "  -14.46%  "
$ws.Range("E37").Value = "  -1.12%  "
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("E39").Value = "  +3.19%  "
$ws.Range("E40").Value = "  -7.00%  "
$ws.Range("E41").Value = "  +1.29%  "
$ws.Range("D42").Value = "'69.32"
$ws.Range("E42").Value = "  -2.96%  "
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("E44").Value = "  -0.99%  "
$ws.Range("D45").Value = "'90.92"
$ws.Range("E45").Value = "  -12.12%  "
$ws.Range("D46").Value = "'12.27"
$ws.Range("E46").Value = "  +1.03%  "
$ws.Range("D47").Value = "'79.87"
$ws.Range("E47").Value = "  +2.21%  "
$ws.Range("D48").Value = "'112.53"
$ws.Range("E48").Value = "  -3.30%  "
$ws.Range("D49").Value = "'8.92"
$ws.Range("E49").Value = "  -2.65%  "
$ws.Range("E50").Value = "  -2.59%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "1.588.74"
$ws.Range("E51").Value = "  +2.24%  "
